$d = $word.ActiveDocument

$replacements = @(
    @("2026-01-15 Thursday", "2026-01-16 Friday"),
    @("743×8=5944", "426×5=2130"),
    @("242×2=484", "541×4=2164"),
    @("591×2=1182", "575×7=4025"),
    @("344×3=1032", "832×9=7488"),
    @("607×6=3642", "924×2=1848"),
    @("183×6=1098", "211×7=1477"),
    @("629×2=1258", "843×6=5058"),
    @("610×5=3050", "992×2=1984"),
    @("889×7=6223", "880×8=7040"),
    @("976×9=8784", "456×9=4104"),
    @("625×7=4375", "631×7=4417"),
    @("365×5=1825", "820×4=3280"),
    @("956×8=7648", "239×6=1434"),
    @("220×3=660", "502×6=3012"),
    @("119×9=1071", "762×9=6858"),
    @("819×8=6552", "416×2=832"),
    @("766×3=2298", "221×6=1326"),
    @("466×7=3262", "533×9=4797"),
    @("602×2=1204", "421×4=1684"),
    @("125×3=375", "545×8=4360"),
    @("673×3=2019", "142×4=568"),
    @("277×7=1939", "304×7=2128"),
    @("285×9=2565", "364×5=1820"),
    @("151×7=1057", "864×6=5184"),
    @("245×5=1225", "645×7=4515")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
